$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A20").Value = "Post 5.jpg"
$ws.Range("B20").Value = 5
$ws.Range("C20").Value = "👫🏻 **Celebrate the unbreakable bond of sibling love**`n⚡ 50%-80% off`n👉🏻 amzaff.to/93m1fy6`n🎁 **Person"
$ws.Range("D20").Value = "Image + Text sent"
$ws.Range("E20").Value = "2025-07-21 12:18:04"

$ws.Range("A21").Value = "Post 3.jpg"
$ws.Range("B21").Value = 3
$ws.Range("C21").Value = "🧶 **amazonkarigar | Threads of heritage**`n👉🏻 amzaff.to/otM7QCn`n🔥 **Latest Arrivals** 🔥`n🚩 **Indian "
$ws.Range("D21").Value = "Image + Text sent"
$ws.Range("E21").Value = "2025-07-21 12:22:05"

$ws.Range("A22").Value = "Post 1.jpg"
$ws.Range("B22").Value = 1
$ws.Range("C22").Value = "🔥 **BLOCKBUSTER DEALS OF THE DAY** 🔥`n💰 **Prime Everyday Savings**`n⚡ Up to 60% off`n+ Extra 5% off,`n-"
$ws.Range("D22").Value = "Image + Text sent"
$ws.Range("E22").Value = "2025-07-21 12:22:05"

$ws.Range("A23").Value = "Post 4.jpg"
$ws.Range("B23").Value = 4
$ws.Range("C23").Value = "📺 **Smart TVs & Projectors**`n⚡ Starting at ₹8,999`n👉🏻 amzaff.to/Kkk481A`n⏱️ **Limited Time Mega Deals"
$ws.Range("D23").Value = "Image + Text sent"
$ws.Range("E23").Value = "2025-07-21 12:22:06"

$ws.Range("A24").Value = "Post 5.jpg"
$ws.Range("B24").Value = 5
$ws.Range("C24").Value = "👫🏻 **Celebrate the unbreakable bond of sibling love**`n⚡ 50%-80% off`n👉🏻 amzaff.to/93m1fy6`n🎁 **Person"
$ws.Range("D24").Value = "Image + Text sent"
$ws.Range("E24").Value = "2025-07-21 12:23:05"

$ws.Range("A25").Value = "Post 2.jpg"
$ws.Range("B25").Value = 2
$ws.Range("C25").Value = "💥 **Mega Electronic Days**`n📣 Sale live till 20th July`n💻 **Laptops, Smartwatches, Headphones & more*"
$ws.Range("D25").Value = "Image + Text sent"
$ws.Range("E25").Value = "2025-07-21 12:24:12"

$ws.Range("A26").Value = "Post 1.jpg"
$ws.Range("B26").Value = 1
$ws.Range("C26").Value = "🔥 **BLOCKBUSTER DEALS OF THE DAY** 🔥`n💰 **Prime Everyday Savings**`n⚡ Up to 60% off`n+ Extra 5% off,`n-"
$ws.Range("D26").Value = "Image + Text sent"
$ws.Range("E26").Value = "2025-07-21 12:24:14"

$ws.Range("A27").Value = "Post 3.jpg"
$ws.Range("B27").Value = 3
$ws.Range("C27").Value = "🧶 **amazonkarigar | Threads of heritage**`n👉🏻 amzaff.to/otM7QCn`n🔥 **Latest Arrivals** 🔥`n🚩 **Indian "
$ws.Range("D27").Value = "Image + Text sent"
$ws.Range("E27").Value = "2025-07-21 12:24:14"

$ws.Range("A28").Value = "Post 4.jpg"
$ws.Range("B28").Value = 4
$ws.Range("C28").Value = "📺 **Smart TVs & Projectors**`n⚡ Starting at ₹8,999`n👉🏻 amzaff.to/Kkk481A`n⏱️ **Limited Time Mega Deals"
$ws.Range("D28").Value = "Image + Text sent"
$ws.Range("E28").Value = "2025-07-21 12:24:14"

$ws.Range("A29").Value = "Post 5.jpg"
$ws.Range("B29").Value = 5
$ws.Range("C29").Value = "👫🏻 **Celebrate the unbreakable bond of sibling love**`n⚡ 50%-80% off`n👉🏻 amzaff.to/93m1fy6`n🎁 **Person"
$ws.Range("D29").Value = "Image + Text sent"
$ws.Range("E29").Value = "2025-07-21 12:24:14"

$ws.Range("A30").Value = "Post 2.jpg"
$ws.Range("B30").Value = 2
$ws.Range("C30").Value = "💥 **Mega Electronic Days**`n📣 Sale live till 20th July`n💻 **Laptops, Smartwatches, Headphones & more*"
$ws.Range("D30").Value = "Image + Text sent"
$ws.Range("E30").Value = "2025-07-21 12:32:25"

$ws.Range("A31").Value = "Post 3.jpg"
$ws.Range("B31").Value = 3
$ws.Range("C31").Value = "🧶 **amazonkarigar | Threads of heritage**`n👉🏻 amzaff.to/otM7QCn`n🔥 **Latest Arrivals** 🔥`n🚩 **Indian "
$ws.Range("D31").Value = "Image + Text sent"
$ws.Range("E31").Value = "2025-07-21 12:36:03"

$ws.Range("A32").Value = "Post 2.jpg"
$ws.Range("B32").Value = 2
$ws.Range("C32").Value = "💥 **Mega Electronic Days**`n📣 Sale live till 20th July`n💻 **Laptops, Smartwatches, Headphones & more*"
$ws.Range("D32").Value = "Image + Text sent"
$ws.Range("E32").Value = "2025-07-21 12:38:08"

$ws.Range("A33").Value = "Post 1.jpg"
$ws.Range("B33").Value = 1
$ws.Range("C33").Value = "🔥 **BLOCKBUSTER DEALS OF THE DAY** 🔥`n💰 **Prime Everyday Savings**`n⚡ Up to 60% off`n+ Extra 5% off,`n-"
$ws.Range("D33").Value = "Image + Text sent"
$ws.Range("E33").Value = "2025-07-21 12:38:09"

$ws.Range("A34").Value = "Post 4.jpg"
$ws.Range("B34").Value = 4
$ws.Range("C34").Value = "📺 **Smart TVs & Projectors**`n⚡ Starting at ₹8,999`n👉🏻 amzaff.to/Kkk481A`n⏱️ **Limited Time Mega Deals"
$ws.Range("D34").Value = "Image + Text sent"
$ws.Range("E34").Value = "2025-07-21 12:38:09"

$ws.Range("A35").Value = "Post 3.jpg"
$ws.Range("B35").Value = 3
$ws.Range("C35").Value = "🧶 **amazonkarigar | Threads of heritage**`n👉🏻 amzaff.to/otM7QCn`n🔥 **Latest Arrivals** 🔥`n🚩 **Indian "
$ws.Range("D35").Value = "Image + Text sent"
$ws.Range("E35").Value = "2025-07-21 12:38:10"

$ws.Range("A36").Value = "Post 5.jpg"
$ws.Range("B36").Value = 5
$ws.Range("C36").Value = "👫🏻 **Celebrate the unbreakable bond of sibling love**`n⚡ 50%-80% off`n👉🏻 amzaff.to/93m1fy6`n🎁 **Person"
$ws.Range("D36").Value = "Image + Text sent"
$ws.Range("E36").Value = "2025-07-21 12:38:10"

$ws.Range("A37").Value = "Post 2.jpg"
$ws.Range("B37").Value = 2
$ws.Range("C37").Value = "💥 **Mega Electronic Days**`n📣 Sale live till 20th July`n💻 **Laptops, Smartwatches, Headphones & more*"
$ws.Range("D37").Value = "Image + Text sent"
$ws.Range("E37").Value = "2025-07-21 12:40:09"

$ws.Range("A38").Value = "Post 1.jpg"
$ws.Range("B38").Value = 1
$ws.Range("C38").Value = "🔥 **BLOCKBUSTER DEALS OF THE DAY** 🔥`n💰 **Prime Everyday Savings**`n⚡ Up to 60% off`n+ Extra 5% off,`n-"
$ws.Range("D38").Value = "Image + Text sent"
$ws.Range("E38").Value = "2025-07-21 12:40:10"

$ws.Range("A39").Value = "Post 3.jpg"
$ws.Range("B39").Value = 3
$ws.Range("C39").Value = "🧶 **amazonkarigar | Threads of heritage**`n👉🏻 amzaff.to/otM7QCn`n🔥 **Latest Arrivals** 🔥`n🚩 **Indian "
$ws.Range("D39").Value = "Image + Text sent"
$ws.Range("E39").Value = "2025-07-21 12:40:10"

$ws.Range("A40").Value = "Post 4.jpg"
$ws.Range("B40").Value = 4
$ws.Range("C40").Value = "📺 **Smart TVs & Projectors**`n⚡ Starting at ₹8,999`n👉🏻 amzaff.to/Kkk481A`n⏱️ **Limited Time Mega Deals"
$ws.Range("D40").Value = "Image + Text sent"
$ws.Range("E40").Value = "2025-07-21 12:40:10"

$ws.Range("A41").Value = "Post 5.jpg"
$ws.Range("B41").Value = 5
$ws.Range("C41").Value = "👫🏻 **Celebrate the unbreakable bond of sibling love**`n⚡ 50%-80% off`n👉🏻 amzaff.to/93m1fy6`n🎁 **Person"
$ws.Range("D41").Value = "Image + Text sent"
$ws.Range("E41").Value = "2025-07-21 12:40:10"

Write-Output "done"
